$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()

$ws.Range("H53").Value = 687.5
$ws.Range("I53").Value = 900
$ws.Range("J53").Value = 616.6667
$ws.Range("K53").Value = 900
$ws.Range("L53").Value = 616.6667
$ws.Range("M53").Value = -263
$ws.Range("N53").Value = -1890.6667

$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("N55").ClearContents()

$ws.Range("H58").Value = 140
$ws.Range("I58").Value = 140
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 420
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -270
$ws.Range("N58").ClearContents()

$ws.Range("H62").Value = 4575
$ws.Range("I62").Value = 5033.3335
$ws.Range("J62").Value = 3887.5
$ws.Range("K62").Value = 5033.3335
$ws.Range("L62").Value = 3887.5
$ws.Range("M62").Value = -4409.3335
$ws.Range("N62").Value = -5135.5

$ws.Range("H65").Value = 4575
$ws.Range("I65").Value = 5033.3335
$ws.Range("J65").Value = 3887.5
$ws.Range("K65").Value = 25166.6675
$ws.Range("L65").Value = 19437.5
$ws.Range("M65").Value = -22046.6675
$ws.Range("N65").Value = -25677.5

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H137").Value = 3834

$ws.Range("H138").Value = 5266968
$ws.Range("I138").Value = 25003688
$ws.Range("J138").Value = 3842.8
$ws.Range("K138").Value = 75011064
$ws.Range("L138").Value = 11528.4
$ws.Range("M138").Value = -75005924
$ws.Range("N138").Value = -21808.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3117.7058
$ws.Range("I132").Value = 853.9231
$ws.Range("K132").Value = 2561.7693
$ws.Range("M132").Value = -31.76929999999993

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3006
$ws.Range("I86").Value = 3006
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 3006
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -1883
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 3006
$ws.Range("I89").Value = 3006
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 15030
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -9414
$ws.Range("N89").ClearContents()

$ws.Range("H99").Value = 1129
$ws.Range("I99").Value = 1129
$ws.Range("K99").Value = 1129
$ws.Range("M99").Value = 369

$ws.Range("H105").Value = 2091.5
$ws.Range("I105").Value = 2033.1428
$ws.Range("K105").Value = 2033.1428
$ws.Range("M105").Value = -286.1428000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 559.6
$ws.Range("I12").Value = 404
$ws.Range("J12").Value = 598.5
$ws.Range("K12").Value = 404
$ws.Range("L12").Value = 598.5
$ws.Range("M12").Value = -234
$ws.Range("N12").Value = -938.5

$ws.Range("H35").Value = 4209.8
$ws.Range("I35").Value = 1025
$ws.Range("K35").Value = 1025
$ws.Range("M35").Value = -731

$ws.Range("H38").Value = 10000
$ws.Range("I38").Value = 10000
$ws.Range("K38").Value = 10000
$ws.Range("M38").Value = -9623

$ws.Range("H46").Value = 10000
$ws.Range("I46").Value = 10000
$ws.Range("K46").Value = 10000
$ws.Range("M46").Value = -9789

$ws.Range("H106").Value = 25750
$ws.Range("J106").Value = 25750
$ws.Range("L106").Value = 25750
$ws.Range("N106").Value = -28274

$ws.Range("H122").Value = 1936.8572
$ws.Range("I122").Value = 1843.1666
$ws.Range("K122").Value = 5529.4998
$ws.Range("M122").Value = -3079.4998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 4597
$ws.Range("I104").Value = 411
$ws.Range("K104").Value = 1233
$ws.Range("M104").Value = 1388

$ws.Range("H117").Value = 2694.8462
$ws.Range("J117").Value = 2502.75
$ws.Range("L117").Value = 7508.25
$ws.Range("N117").Value = -14392.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 21499.5
$ws.Range("J53").Value = 21499.5
$ws.Range("L53").Value = 21499.5
$ws.Range("N53").Value = -22761.5

$ws.Range("H57").Value = 30030
$ws.Range("J57").Value = 30030
$ws.Range("L57").Value = 30030
$ws.Range("N57").Value = -31670

$ws.Range("H80").Value = 3985.6667
$ws.Range("I80").Value = 3958
$ws.Range("J80").Value = 3999.5
$ws.Range("K80").Value = 3958
$ws.Range("L80").Value = 3999.5
$ws.Range("M80").Value = -2960
$ws.Range("N80").Value = -5995.5

$ws.Range("H83").Value = 3985.6667
$ws.Range("I83").Value = 3958
$ws.Range("J83").Value = 3999.5
$ws.Range("K83").Value = 19790
$ws.Range("L83").Value = 19997.5
$ws.Range("M83").Value = -14798
$ws.Range("N83").Value = -29981.5

$ws.Range("H97").Value = 798.3333
$ws.Range("I97").Value = 699.5
$ws.Range("K97").Value = 699.5
$ws.Range("M97").Value = -203.5

$ws.Range("H132").Value = 3588.8823
$ws.Range("I132").Value = 1811.1818
$ws.Range("J132").Value = 6848
$ws.Range("K132").Value = 5433.5454
$ws.Range("L132").Value = 20544
$ws.Range("M132").Value = -2903.5454
$ws.Range("N132").Value = -25604

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()

$ws.Range("H61").Value = 1457.3334
$ws.Range("I61").Value = 1221.2222
$ws.Range("K61").Value = 1221.2222
$ws.Range("M61").Value = -1019.2222

$ws.Range("H68").Value = 3135.111
$ws.Range("I68").Value = 3145.1428
$ws.Range("J68").Value = 3100
$ws.Range("K68").Value = 3145.1428
$ws.Range("L68").Value = 3100
$ws.Range("M68").Value = -2396.1428
$ws.Range("N68").Value = -4598

$ws.Range("H71").Value = 3135.111
$ws.Range("I71").Value = 3145.1428
$ws.Range("J71").Value = 3100
$ws.Range("K71").Value = 15725.714
$ws.Range("L71").Value = 15500
$ws.Range("M71").Value = -11981.714
$ws.Range("N71").Value = -22988

$ws.Range("H113").Value = 1457.3334
$ws.Range("I113").Value = 1221.2222
$ws.Range("K113").Value = 1221.2222
$ws.Range("M113").Value = 948.7778000000001

$ws.Range("H132").Value = 7861.5386
$ws.Range("I132").Value = 5619.5
$ws.Range("K132").Value = 16858.5
$ws.Range("M132").Value = -14328.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 6000
$ws.Range("J12").Value = 7000
$ws.Range("L12").Value = 7000
$ws.Range("N12").Value = -7284

$ws.Range("H122").Value = 1999
$ws.Range("I122").Value = 1998.5
$ws.Range("K122").Value = 5995.5
$ws.Range("M122").Value = -3545.5
